$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").Value = -7.6945
$ws.Range("B10").Value = 5.6295
$ws.Range("B12").Value = 5.291799999999999
$ws.Range("D12").Value = -5.761099999999998
$ws.Range("D17").Value = -8.365899999999993
$ws.Range("B18").Value = 6.928799999999993
$ws.Range("D26").Value = -7.215400000000008
$ws.Range("D27").Value = -8.101699999999999
$ws.Range("D28").Value = -8.301299999999998
$ws.Range("B37").Value = 8.604400000000005
$ws.Range("D37").Value = -8.113400000000002
$ws.Range("B55").Value = 6.508999999999995
$ws.Range("D65").Value = -7.944700000000003
$ws.Range("B68").Value = 4.935099999999998
$ws.Range("D73").Value = -8.399399999999995
$ws.Range("B77").Value = 9.1511
$ws.Range("B78").Value = 9.507399999999997
$ws.Range("D84").Value = -8.019500000000004
$ws.Range("D85").Value = -8.961199999999995
$ws.Range("D93").Value = -6.799899999999993
$ws.Range("D95").Value = -7.431100000000002
$ws.Range("D98").Value = -7.209500000000003
$ws.Range("D99").Value = -7.932000000000004
$ws.Range("D101").Value = -7.810299999999996

